# Applies the "fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về
# chi tiêu" edit across all affected sheets of the monthly SÓC TRĂNG report.
#
# Helper: write a full row of a worksheet from a positional array. $null
# entries clear the cell; $textCols marks 1-based columns whose
# NumberFormat must be forced to "@" (text) first so date-look-alike
# strings ("08-04-2024", "2024-05-25", ...) are not auto-converted into
# serial date numbers by Excel.
function Set-RowValues {
    param($ws, $row, $values, $textCols = @())
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $v = $values[$i]
        if ($null -eq $v) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet 1: CHI TIẾT DOANH THU — 6 new order rows (625,626,627,628,636,637)
# inserted before the totals row, which moves from row 4 to row 10.
# ======================================================================
$ws1 = $wb.Worksheets.Item(1)
$sheet1Rows = @{
    4 = @("08-04-2024", "HD-LUXURY", 625, "SÓC TRĂNG", "Cắt mí", "nguyễn thị mỹ chăm", "Cá nhân", "Lâm Thị Mỹ Hằng", $null, "Lê Hoàng Thanh", 6000000, 6000000, "Lâm Thị Mỹ Hằng", $null, 6000000, 0, 6000000, 0, "Kha Như Huỳnh ", $null, 50000, 0)
    5 = @("08-04-2024", "HD-LUXURY", 626, "SÓC TRĂNG", "nhấn đồng tiền", "nguyễn thị mỹ trinh", "Cá nhân", "Lâm Thị Mỹ Hằng", 7000000, $null, $null, 7000000, "Nguyễn Hoàng Yến Quyên", $null, 7000000, 0, 7000000, 0, "Trần Khánh Hiệp", $null, $null, $null)
    6 = @("08-04-2024", "HD-LUXURY", 627, "SÓC TRĂNG", "Cắt mí", "tạ duy hoàng ", "Cá nhân", "Lê Đình Hậu", 6000000, $null, $null, 6000000, "Lâm Thị Mỹ Hằng", $null, 6000000, 0, 6000000, 0, $null, "Kha Như Huỳnh ", 50000, 0)
    7 = @("08-04-2024", "HD-LUXURY", 628, "SÓC TRĂNG", "Cắt mí", "nguyễn thị lệ trang", "Cá nhân", "Lê Đình Hậu", 4000000, $null, $null, 4000000, "Nguyễn Hoàng Yến Quyên", $null, 4000000, 0, 4000000, 0, "Trần Khánh Hiệp", $null, 50000, 0)
    8 = @("08-09-2024", "HD-LUXURY", 636, "SÓC TRĂNG", "cấy mỡ mặt ", "thạch thị siêu", "CTV", "Thạch Hoàng Nhân", 36000000, $null, $null, 36000000, $null, $null, 36000000, 0, 36000000, 0, $null, $null, $null, $null)
    9 = @("08-09-2024", "HD-LUXURY", 637, "SÓC TRĂNG", "Nâng mũi", "thạch thị siêu", "CTV", "Thạch Hoàng Nhân", 40000000, $null, $null, 40000000, $null, $null, 14000000, 0, 14000000, 26000000, $null, $null, 100000, 50000)
    10 = @($null, "Tổng", 8, $null, $null, $null, $null, $null, 107000000, $null, 6000000, 113000000, $null, $null, 85000000, 0, 85000000, 28000000, $null, $null, 400000, 50000)
}
$sheet1TextCols = @(1)
foreach ($r in ($sheet1Rows.Keys | Sort-Object)) {
    Set-RowValues $ws1 $r $sheet1Rows[$r] $sheet1TextCols
}

# ======================================================================
# Sheet 2: CHI TIẾT VỀ THU NỢ — 5 new debt-collection rows; totals move
# from row 2 to row 7.
# ======================================================================
$ws2 = $wb.Worksheets.Item(2)
$sheet2Rows = @{
    2 = @("TN", 176, "08-04-2024", "SÓC TRĂNG", "HD-LUXURY-428", "2024-05-25", 500000)
    3 = @("TN", 177, "08-05-2024", "SÓC TRĂNG", "HD-LUXURY-611", "2024-07-31", 8000000)
    4 = @("TN", 178, "08-06-2024", "SÓC TRĂNG", "HD-LUXURY-356", "2024-04-29", 2000000)
    5 = @("TN", 179, "08-09-2024", "SÓC TRĂNG", "HD-LUXURY-500", "2024-06-30", $null)
    6 = @("TN", 180, "08-09-2024", "SÓC TRĂNG", "HD-LUXURY-500", "2024-06-30", 2000000)
    7 = @("Tổng", 5, $null, $null, $null, $null, 12500000)
}
$sheet2TextCols = @(3, 6)
foreach ($r in ($sheet2Rows.Keys | Sort-Object)) {
    Set-RowValues $ws2 $r $sheet2Rows[$r] $sheet2TextCols
}

# ======================================================================
# Sheet 3: CHI TIẾT CHI TIÊU — new "Ghi chú" column (G) + 8 new expense
# rows; totals move from row 6 to row 14.
# ======================================================================
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(1, 7).Value = "Ghi chú"
$ws3.Cells.Item(2, 7).ClearContents()
$ws3.Cells.Item(3, 7).ClearContents()
$ws3.Cells.Item(4, 7).ClearContents()
$ws3.Cells.Item(5, 7).ClearContents()

$sheet3Rows = @{
    6 = @("CT", 755, "08-04-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 940000, $null)
    7 = @("CT", 756, "08-05-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 500000, $null)
    8 = @("CT", 761, "08-06-2024", "SÓC TRĂNG", "Chi Phí CTV", 400000, $null)
    9 = @("CT", 766, "08-08-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 400000, $null)
    10 = @("CT", 767, "08-08-2024", "SÓC TRĂNG", "Chi Phí CTV", 800000, $null)
    11 = @("CT", 770, "08-09-2024", "SÓC TRĂNG", "Chi Phí Vận Hành", 4700000, $null)
    12 = @("CT", 771, "08-09-2024", "SÓC TRĂNG", "Chi Phí Vận Hành", 1800000, $null)
    13 = @("CT", 772, "08-09-2024", "SÓC TRĂNG", "Chi Phí CTV", 10400000, "triết khấu tiền phần trăm khách và khách tn")
    14 = @("Tổng", 12, $null, $null, $null, 21500000, $null)
}
$sheet3TextCols = @(3)
foreach ($r in ($sheet3Rows.Keys | Sort-Object)) {
    Set-RowValues $ws3 $r $sheet3Rows[$r] $sheet3TextCols
}

# ======================================================================
# Sheet 4: DOANH SỐ CÁ NHÂN — per-employee figures recomputed for the
# fuller month; roster also changes (one fewer row than before), so the
# old trailing row (11) is cleared after the new rows are written.
# ======================================================================
$ws4 = $wb.Worksheets.Item(4)
$sheet4Rows = @{
    2 = @("Cô Siêng giúp Việc", 0, 0, 0, 0, 0, 0, 0, 0, 500000)
    3 = @("Kha Như Huỳnh ", 0, 0, 0, 0, 3, 200000, 1, 0, 0)
    4 = @("Lâm Thị Mỹ Hằng", 7000000, 0, 18000000, 0, 0, 0, 0, 0, 0)
    5 = @("Lê Hoàng Thanh", 0, 6000000, 0, 0, 0, 0, 0, 0, 0)
    6 = @("Lê Đình Hậu", 24000000, 0, 0, 0, 0, 0, 0, 0, 0)
    7 = @("Nguyễn Hoàng Yến Quyên", 0, 0, 17000000, 0, 0, 0, 0, 0, 0)
    8 = @("Thạch Hoàng Nhân", 76000000, 0, 0, 0, 0, 0, 0, 0, 12000000)
    9 = @("Trần Khánh Hiệp", 0, 0, 0, 0, 2, 50000, 0, 0, 0)
    10 = @("Tổng", 107000000, 6000000, 35000000, 0, 5, 250000, 1, 0, 12500000)
}
foreach ($r in ($sheet4Rows.Keys | Sort-Object)) {
    Set-RowValues $ws4 $r $sheet4Rows[$r]
}
$ws4.Rows.Item(11).ClearContents()

# ======================================================================
# Sheet 5: CHI TIÊU TỔNG HỢP — expense-by-category rollup gains a "Chi
# Phí CTV" category row; totals move from row 6 to row 7.
# ======================================================================
$ws5 = $wb.Worksheets.Item(5)
$sheet5Rows = @{
    2 = @("Chi Phí CTV", 11600000)
    3 = @("Chi Phí Sinh Hoạt Tại Cơ Sở", 2200000)
    4 = @("Chi Phí Vận Hành", 7000000)
    5 = @("Trang thiết bị Y Tế", 700000)
    6 = @("Blank", 0)
    7 = @("Tổng cộng", 21500000)
}
foreach ($r in ($sheet5Rows.Keys | Sort-Object)) {
    Set-RowValues $ws5 $r $sheet5Rows[$r]
}

# ======================================================================
# Sheet 6: LŨY KẾ NGÀY — daily running totals extended through 08-09;
# rows 2-3 (08-02, 08-03) are untouched, new rows 4-8 are inserted before
# the totals row, which moves from row 4 to row 9.
# ======================================================================
$ws6 = $wb.Worksheets.Item(6)
$sheet6Rows = @{
    4 = @("08-04-2024", 23000000, 23000000, 4, 500000, 940000, 22560000)
    5 = @("08-05-2024", 0, 0, 0, 8000000, 500000, 7500000)
    6 = @("08-06-2024", 0, 0, 0, 2000000, 400000, 1600000)
    7 = @("08-08-2024", 0, 0, 0, 0, 1200000, -1200000)
    8 = @("08-09-2024", 76000000, 50000000, 2, 2000000, 16900000, 35100000)
    9 = @("Tổng", 113000000, 85000000, 8, 12500000, 21500000, 76000000)
}
$sheet6TextCols = @(1)
foreach ($r in ($sheet6Rows.Keys | Sort-Object)) {
    Set-RowValues $ws6 $r $sheet6Rows[$r] $sheet6TextCols
}

# ======================================================================
# Sheet 7: QUỸ LƯƠNG — payroll fund re-allocated across the full month;
# roster/row layout is unchanged, only the "Tổng lương tại SÓC TRĂNG"
# column (C) values move.
# ======================================================================
$ws7 = $wb.Worksheets.Item(7)
$ws7.Cells.Item(4, 3).Value = 857142.8571428573
$ws7.Cells.Item(8, 3).Value = 2557142.857142857
$ws7.Cells.Item(9, 3).Value = 2752380.952380952
$ws7.Cells.Item(10, 3).Value = 1628571.428571429
$ws7.Cells.Item(11, 3).Value = 4285714.285714285
$ws7.Cells.Item(17, 3).Value = 5571428.571428571
$ws7.Cells.Item(18, 3).Value = 1472857.142857143
$ws7.Cells.Item(19, 3).Value = 1117142.857142857
$ws7.Cells.Item(20, 3).Value = 1300000
$ws7.Cells.Item(21, 3).Value = 852857.1428571428
$ws7.Cells.Item(22, 3).Value = 22395238.09523809

# ======================================================================
# Sheet 8: LỢI NHUẬN — base-level profitability summary recomputed for
# the extended month; single data row (2), layout unchanged.
# ======================================================================
$ws8 = $wb.Worksheets.Item(8)
$ws8.Cells.Item(2, 2).Value = 113000000
$ws8.Cells.Item(2, 3).Value = 85000000
$ws8.Cells.Item(2, 4).Value = 0.7522123893805309
$ws8.Cells.Item(2, 5).Value = 0.2477876106194691
$ws8.Cells.Item(2, 6).Value = 12500000
$ws8.Cells.Item(2, 7).Value = 97500000
$ws8.Cells.Item(2, 8).Value = 21500000
$ws8.Cells.Item(2, 9).Value = 22395238.09523809
$ws8.Cells.Item(2, 10).Value = 43895238.09523809
$ws8.Cells.Item(2, 11).Value = 53604761.90476191
$ws8.Cells.Item(2, 12).Value = 0.5497924297924298

